$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CreateNewRTIEmployees")
$ws.Activate()

# Update ContractualHours (N2): numeric 40 -> text "16"
$ws.Range("N2").Value = "16"

# Update the Email/UserName shared value (D2 & E2 both point at the same
# shared string): "rtiemployee1003@xcdgmail.com" -> "rtiemployee1014@xcdgmail.com"
$ws.Range("D2").Value = "rtiemployee1014@xcdgmail.com"
$ws.Range("E2").Value = "rtiemployee1014@xcdgmail.com"

# Update LastName (C2): "RTI employee 1002" -> "RTI employee 1014"
$ws.Range("C2").Value = "RTI employee 1014"

# Update view: select D2 (was R5) and scroll the window back to A1 (was topLeftCell R1)
$ws.Range("A1").Select()
$ws.Range("D2").Select()
